$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "Complaint Number" -> "Account ID"
$ws.Range("A1").Value = "Account ID"

# "Status" header stays as-is in B1

# New "Remarks" header in C1, matching the bold style already used by A1/B1
$ws.Range("C1").Value = "Remarks"
$ws.Range("C1").Font.Bold = $true

# Shrink the saved window width
$excel.ActiveWindow.Width = 18350
